# Elimna EC anteriores y se agregan nuevos, se modifica base de datos
#
# Updates the "Periodo Mora" (E16:E21) / "Valor Mora" (F16:F21) block on the
# account-statement sheet: the previous payment-period rows (2110..2203,
# oldest first) are replaced by the new period list in reverse
# (2203..2110, newest first), carrying the corresponding arrears amount
# along with the period that used to sit on the last row.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$periodos = @("2203", "2202", "2201", "2112", "2111", "2110")
$valores  = @(32707, 36341, 36341, 36341, 36341, 36341)

for ($i = 0; $i -lt $periodos.Length; $i++) {
    $row = 16 + $i
    $ws.Cells.Item($row, 5).Value = $periodos[$i]
    $ws.Cells.Item($row, 6).Value = $valores[$i]
}
